$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - new "result" summary cells per site
$ws.Range("A2").Value = "모두 비허용"
$ws.Range("D2").Value = "존재x"
$ws.Range("G2").Value = "모두 허용"
$ws.Range("J2").Value = "존재x"
$ws.Range("M2").Value = "존재 및 모두 허용"

# Row 3 - supporting robots.txt excerpt lines
$ws.Range("A3").Value = "User-agent: *"
$ws.Range("G3").Value = "User-agent: * "

# Row 4 - supporting robots.txt excerpt lines
$ws.Range("A4").Value = "Disallow: /"
$ws.Range("G4").Value = "Allow : /"

# Restore the author's last selection
$ws.Range("M6").Select()
